$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1749544.4
$ws.Range("I70").Value = 2842050.5
$ws.Range("J70").Value = 1534.4
$ws.Range("K70").Value = 8526151.5
$ws.Range("L70").Value = 4603.200000000001
$ws.Range("M70").Value = -8525881.5
$ws.Range("N70").Value = -5143.200000000001

$ws.Range("H73").Value = 1749544.4
$ws.Range("I73").Value = 2842050.5
$ws.Range("J73").Value = 1534.4
$ws.Range("K73").Value = 8526151.5
$ws.Range("L73").Value = 4603.200000000001
$ws.Range("M73").Value = -8525215.5
$ws.Range("N73").Value = -6475.200000000001

$ws.Range("H116").Value = 2886.5386
$ws.Range("I116").Value = 2787.0527
$ws.Range("J116").Value = 3156.5715
$ws.Range("K116").Value = 2787.0527
$ws.Range("L116").Value = 3156.5715
$ws.Range("M116").Value = 654.9472999999998
$ws.Range("N116").Value = -10040.5715

$ws.Range("H132").Value = 4236.5312
$ws.Range("I132").Value = 4169.033
$ws.Range("J132").Value = 5249
$ws.Range("K132").Value = 12507.099
$ws.Range("L132").Value = 15747
$ws.Range("M132").Value = -9977.099000000002
$ws.Range("N132").Value = -20807

$ws.Range("H136").Value = 52000
$ws.Range("J136").Value = 52000
$ws.Range("L136").Value = 52000
$ws.Range("N136").Value = -62200

$ws.Range("H138").Value = 169201.14
$ws.Range("I138").Value = 2175.9
$ws.Range("J138").Value = 243434.58
$ws.Range("K138").Value = 6527.700000000001
$ws.Range("L138").Value = 730303.74
$ws.Range("M138").Value = -1387.700000000001
$ws.Range("N138").Value = -740583.74

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 40980
$ws.Range("J7").Value = 40980
$ws.Range("L7").Value = 40980
$ws.Range("N7").Value = -41208

$ws.Range("H61").Value = 3577.8
$ws.Range("I61").Value = 3577.8
$ws.Range("K61").Value = 3577.8
$ws.Range("M61").Value = -3365.8

$ws.Range("H63").Value = 5798.643
$ws.Range("I63").Value = 3215.8333
$ws.Range("J63").Value = 7735.75
$ws.Range("K63").Value = 3215.8333
$ws.Range("L63").Value = 7735.75
$ws.Range("M63").Value = -2529.8333
$ws.Range("N63").Value = -9107.75

$ws.Range("H66").Value = 5798.643
$ws.Range("I66").Value = 3215.8333
$ws.Range("J66").Value = 7735.75
$ws.Range("K66").Value = 16079.1665
$ws.Range("L66").Value = 38678.75
$ws.Range("M66").Value = -12647.1665
$ws.Range("N66").Value = -45542.75

$ws.Range("H69").Value = 70000
$ws.Range("J69").Value = 70000
$ws.Range("L69").Value = 70000
$ws.Range("N69").Value = -71498

$ws.Range("H72").Value = 70000
$ws.Range("J72").Value = 70000
$ws.Range("L72").Value = 210000
$ws.Range("N72").Value = -217488

$ws.Range("H93").Value = 72965.336
$ws.Range("J93").Value = 72965.336
$ws.Range("L93").Value = 72965.336
$ws.Range("N93").Value = -77957.336

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""

$ws.Range("H113").Value = 39999
$ws.Range("J113").Value = 39999
$ws.Range("L113").Value = 39999
$ws.Range("N113").Value = -48677

$ws.Range("H136").Value = 3577.8
$ws.Range("I136").Value = 3577.8
$ws.Range("K136").Value = 10733.4
$ws.Range("M136").Value = -8183.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 37000
$ws.Range("J2").Value = 37000
$ws.Range("L2").Value = 37000
$ws.Range("N2").Value = -37226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1200.8889
$ws.Range("I58").Value = 815.5
$ws.Range("J58").Value = 1311
$ws.Range("K58").Value = 815.5
$ws.Range("L58").Value = 1311
$ws.Range("M58").Value = -612.5
$ws.Range("N58").Value = -1717

$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -4748

$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -23740

$ws.Range("H99").Value = 1913.0435
$ws.Range("I99").Value = 1620
$ws.Range("J99").Value = 1994.4445
$ws.Range("K99").Value = 1620
$ws.Range("L99").Value = 1994.4445
$ws.Range("M99").Value = -122
$ws.Range("N99").Value = -4990.4445

$ws.Range("H122").Value = 1643.52
$ws.Range("I122").Value = 1012.5714
$ws.Range("J122").Value = 1888.8889
$ws.Range("K122").Value = 3037.7142
$ws.Range("L122").Value = 5666.6667
$ws.Range("M122").Value = -587.7142000000003
$ws.Range("N122").Value = -10566.6667

$ws.Range("H126").Value = 1913.0435
$ws.Range("I126").Value = 1620
$ws.Range("J126").Value = 1994.4445
$ws.Range("K126").Value = 4860
$ws.Range("L126").Value = 5983.333500000001
$ws.Range("M126").Value = -2390
$ws.Range("N126").Value = -10923.3335

$ws.Range("H132").Value = 8335180.5
$ws.Range("I132").Value = 991.4
$ws.Range("J132").Value = 16669370
$ws.Range("K132").Value = 2974.2
$ws.Range("L132").Value = 50008110
$ws.Range("M132").Value = -444.1999999999998
$ws.Range("N132").Value = -50013170

$ws.Range("H134").Value = 2438.2856
$ws.Range("I134").Value = 2383.6
$ws.Range("J134").Value = 2575
$ws.Range("K134").Value = 7150.799999999999
$ws.Range("L134").Value = 7725
$ws.Range("M134").Value = -4615.799999999999
$ws.Range("N134").Value = -12795

$ws.Range("H136").Value = 1200.8889
$ws.Range("I136").Value = 815.5
$ws.Range("J136").Value = 1311
$ws.Range("K136").Value = 2446.5
$ws.Range("L136").Value = 3933
$ws.Range("M136").Value = 103.5
$ws.Range("N136").Value = -9033

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 606.5
$ws.Range("I82").Value = 606.5
$ws.Range("K82").Value = 1819.5
$ws.Range("M82").Value = -1413.5

$ws.Range("H85").Value = 606.5
$ws.Range("I85").Value = 606.5
$ws.Range("K85").Value = 1819.5
$ws.Range("M85").Value = -415.5

$ws.Range("H113").Value = 2756.3333
$ws.Range("J113").Value = 3684.5
$ws.Range("L113").Value = 11053.5
$ws.Range("N113").Value = -15393.5

$ws.Range("H134").Value = 3946.25
$ws.Range("I134").Value = 2285.476
$ws.Range("J134").Value = 8928.571
$ws.Range("K134").Value = 6856.428
$ws.Range("L134").Value = 26785.713
$ws.Range("M134").Value = -1786.428
$ws.Range("N134").Value = -36925.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2447.6843
$ws.Range("I132").Value = 1717.6364
$ws.Range("K132").Value = 5152.9092
$ws.Range("M132").Value = -2622.9092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3837.625
$ws.Range("I16").Value = 3900
$ws.Range("J16").Value = 3733.6667
$ws.Range("K16").Value = 3900
$ws.Range("L16").Value = 3733.6667
$ws.Range("M16").Value = -3730
$ws.Range("N16").Value = -4073.6667

$ws.Range("H68").Value = 3827.2727
$ws.Range("I68").Value = 2421.4285
$ws.Range("J68").Value = 4863.1577
$ws.Range("K68").Value = 2421.4285
$ws.Range("L68").Value = 4863.1577
$ws.Range("M68").Value = -1672.4285
$ws.Range("N68").Value = -6361.1577

$ws.Range("H71").Value = 3827.2727
$ws.Range("I71").Value = 2421.4285
$ws.Range("J71").Value = 4863.1577
$ws.Range("K71").Value = 12107.1425
$ws.Range("L71").Value = 24315.7885
$ws.Range("M71").Value = -8363.1425
$ws.Range("N71").Value = -31803.7885

$ws.Range("H132").Value = 3969.4482
$ws.Range("I132").Value = 3189.0625
$ws.Range("J132").Value = 4929.923
$ws.Range("K132").Value = 9567.1875
$ws.Range("L132").Value = 14789.769
$ws.Range("M132").Value = -7037.1875
$ws.Range("N132").Value = -19849.769

$ws.Range("H136").Value = 2932.5
$ws.Range("I136").Value = 3472.1428
$ws.Range("J136").Value = 1673.3334
$ws.Range("K136").Value = 10416.4284
$ws.Range("L136").Value = 5020.0002
$ws.Range("M136").Value = -7866.428400000001
$ws.Range("N136").Value = -10120.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 77600
$ws.Range("I62").Value = 3200
$ws.Range("J62").Value = 152000
$ws.Range("K62").Value = 3200
$ws.Range("L62").Value = 152000
$ws.Range("M62").Value = -2576
$ws.Range("N62").Value = -153248

$ws.Range("H65").Value = 77600
$ws.Range("I65").Value = 3200
$ws.Range("J65").Value = 152000
$ws.Range("K65").Value = 16000
$ws.Range("L65").Value = 760000
$ws.Range("M65").Value = -12880
$ws.Range("N65").Value = -766240

$ws.Range("H70").Value = 79098.336
$ws.Range("J70").Value = 79098.336
$ws.Range("L70").Value = 79098.336
$ws.Range("N70").Value = -79728.336

$ws.Range("H73").Value = 79098.336
$ws.Range("J73").Value = 79098.336
$ws.Range("L73").Value = 79098.336
$ws.Range("N73").Value = -81282.336

$ws.Range("H129").Value = 70214.5
$ws.Range("J129").Value = 70214.5
$ws.Range("L129").Value = 70214.5
$ws.Range("N129").Value = -80214.5

$ws.Range("H136").Value = 3210.6667
$ws.Range("I136").Value = 2951.3125
$ws.Range("J136").Value = 3418.15
$ws.Range("K136").Value = 8853.9375
$ws.Range("L136").Value = 10254.45
$ws.Range("M136").Value = -6303.9375
$ws.Range("N136").Value = -15354.45
